$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Add bookmarkStart "_Hlk71808202" right before the "Student name:" run
#    (start of paragraph 2, collapsed range). Bookmarks.Add on a collapsed
#    range mints both a bookmarkStart and an auto-paired bookmarkEnd at the
#    same spot; the real bookmarkEnd belongs at the tail of the document
#    (added later), so the auto-paired one is removed once the real one is
#    in place.
# ---------------------------------------------------------------------------
$pStudent = $d.Paragraphs(2)
$rBookmarkStart = $d.Range($pStudent.Range.Start, $pStudent.Range.Start)
$d.Bookmarks.Add("_Hlk71808202", $rBookmarkStart) | Out-Null

# ---------------------------------------------------------------------------
# 2) "January" -> "February" in the Month paragraph.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("January", $true, $false, $false, $false, $false, `
    $true, 1, $false, "February", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Remove " (Week 15 - 18)" -> " " after "February 2021"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(" (Week 15 " + [char]8211 + " 18)", $true, $false, $false, $false, $false, `
    $true, 1, $false, " ", 2) | Out-Null

# ---------------------------------------------------------------------------
# 4) Achievements paragraph text.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute( `
    "This month, I focused on learning about Raspberry Pis and macro cameras, mostly from YouTube tutorials and looking at coding for same. ", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "This month I focused on researching the different model creations for how the image would reference. ", `
    2) | Out-Null

# ---------------------------------------------------------------------------
# 5) Reflection paragraph text (spans two runs -> merges into one).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute( `
    "I think I focused too much on looking at said tutorials when really I should have been looking at machine learning. ", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "I ended up deciding to use TensorFlow Lite as this seems the most straightforward and less complicated way to do. ", `
    2) | Out-Null

# ---------------------------------------------------------------------------
# 6) Intended Changes paragraph text.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute( `
    "After speaking with my supervisor, I realised I was focusing on the wrong aspects of the project. For February, I want to get through the machine learning and start properly coding in March. ", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "After talking with my supervisor, Paul said that working with TensorFlow Lite was probably easier than trying to create a model from scratch, which would severely hinder my time. ", `
    2) | Out-Null

Write-Output "done phase 1"

# ---------------------------------------------------------------------------
# 7) Restructure the "Items discussed" / (delete "Paul also asked...") /
#    "Action Items" paragraphs (paragraphs 13-15) into the new 2-paragraph
#    form, including the spell-check proofErr wrapper and the bookmarkEnd.
# ---------------------------------------------------------------------------
$pItems = $d.Paragraphs(13)
$pAction = $d.Paragraphs(15)
$rBlock = $d.Range($pItems.Range.Start, $pAction.Range.End)

$newXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t>I</w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t>tems</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t xml:space="preserve"> discussed:</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t>Paul and I discussed</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t xml:space="preserve"> TensorFlow Lite and sent me articles to look towards when creating the model. </w:t>
  </w:r>
</w:p>
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:r>
    <w:rPr>
      <w:b/>
      <w:bCs/>
    </w:rPr>
    <w:t>Action Items:</w:t>
  </w:r>
  <w:r>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t xml:space="preserve">I am going </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t>over the contents of what Paul sent.</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:bookmarkEnd w:id="0"/>
</w:p>
'@

$rBlock.InsertXML($newXml) | Out-Null

Write-Output "done phase 2"
